$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style/number format) for the new rows 58-64 by copying
# row 57 formats down, matching existing date (col A) / number (cols B:E) styles.
$ws.Range("A57:E57").Copy()
$ws.Range("A58:E64").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update / add the monthly ECB standard-deviation figures.
$ws.Range("C5").Value = 5.647674560546875
$ws.Range("E5").Value = 5.6075396537780762
$ws.Range("C6").Value = 5.9655489921569824
$ws.Range("E6").Value = 5.6186976432800293
$ws.Range("C7").Value = 5.5414590835571289
$ws.Range("E7").Value = 5.5496826171875
$ws.Range("C8").Value = 5.9174699783325195
$ws.Range("E8").Value = 5.538426399230957
$ws.Range("C9").Value = 4.9655461311340332
$ws.Range("E9").Value = 5.6199564933776855
$ws.Range("C10").Value = 5.6744871139526367
$ws.Range("E10").Value = 5.5932760238647461
$ws.Range("C11").Value = 5.1355934143066406
$ws.Range("E11").Value = 5.5838813781738281
$ws.Range("C12").Value = 5.4596328735351562
$ws.Range("E12").Value = 5.6361603736877441
$ws.Range("C13").Value = 6.2721972465515137
$ws.Range("E13").Value = 5.640160083770752
$ws.Range("C14").Value = 5.4075484275817871
$ws.Range("E14").Value = 5.7810730934143066
$ws.Range("C15").Value = 5.8810000419616699
$ws.Range("E15").Value = 5.796140193939209
$ws.Range("C16").Value = 6.0119681358337402
$ws.Range("E16").Value = 5.7908806800842285
$ws.Range("C17").Value = 5.953467845916748
$ws.Range("E17").Value = 5.7143263816833496
$ws.Range("C18").Value = 6.2337617874145508
$ws.Range("E18").Value = 5.581913948059082
$ws.Range("C19").Value = 5.8100934028625488
$ws.Range("E19").Value = 5.5591616630554199
$ws.Range("C20").Value = 5.0882563591003418
$ws.Range("E20").Value = 5.5207943916320801
$ws.Range("C21").Value = 4.7706451416015625
$ws.Range("E21").Value = 5.4479928016662598
$ws.Range("C22").Value = 5.0804829597473145
$ws.Range("E22").Value = 5.400670051574707
$ws.Range("C23").Value = 5.2027802467346191
$ws.Range("E23").Value = 5.2951059341430664
$ws.Range("C24").Value = 5.5356922149658203
$ws.Range("E24").Value = 5.2893743515014648
$ws.Range("C25").Value = 5.3567571640014648
$ws.Range("E25").Value = 5.394986629486084
$ws.Range("C26").Value = 5.5275616645812988
$ws.Range("E26").Value = 5.5882949829101562
$ws.Range("C27").Value = 5.2836828231811523
$ws.Range("E27").Value = 5.688011646270752
$ws.Range("C28").Value = 5.758509635925293
$ws.Range("E28").Value = 5.781522274017334
$ws.Range("C29").Value = 6.0387692451477051
$ws.Range("E29").Value = 5.8913154602050781
$ws.Range("C30").Value = 6.5104174613952637
$ws.Range("E30").Value = 6.0332798957824707
$ws.Range("C31").Value = 5.9779353141784668
$ws.Range("E31").Value = 6.1529388427734375
$ws.Range("C32").Value = 6.044374942779541
$ws.Range("E32").Value = 6.266258716583252
$ws.Range("C33").Value = 6.523831844329834
$ws.Range("E33").Value = 6.3415679931640625
$ws.Range("C34").Value = 6.6344366073608398
$ws.Range("E34").Value = 6.3831486701965332
$ws.Range("C35").Value = 6.6044936180114746
$ws.Range("E35").Value = 6.3860697746276855
$ws.Range("C36").Value = 6.3035612106323242
$ws.Range("E36").Value = 6.4491972923278809
$ws.Range("C37").Value = 6.4362921714782715
$ws.Range("E37").Value = 6.4834651947021484
$ws.Range("C38").Value = 6.4129962921142578
$ws.Range("E38").Value = 6.4698481559753418
$ws.Range("C39").Value = 6.536707878112793
$ws.Range("E39").Value = 6.4328503608703613
$ws.Range("C40").Value = 6.5460801124572754
$ws.Range("E40").Value = 6.3719758987426758
$ws.Range("C41").Value = 6.3527865409851074
$ws.Range("E41").Value = 6.3504338264465332
$ws.Range("C42").Value = 6.4012770652770996
$ws.Range("E42").Value = 6.2584586143493652
$ws.Range("C43").Value = 6.3014569282531738
$ws.Range("E43").Value = 6.1858644485473633
$ws.Range("C44").Value = 6.0566267967224121
$ws.Range("E44").Value = 6.0964183807373047
$ws.Range("C45").Value = 6.1096820831298828
$ws.Range("E45").Value = 5.9534921646118164
$ws.Range("C46").Value = 5.608515739440918
$ws.Range("E46").Value = 5.8527154922485352
$ws.Range("C47").Value = 5.759648323059082
$ws.Range("E47").Value = 5.7476334571838379
$ws.Range("C48").Value = 5.7316904067993164
$ws.Range("E48").Value = 5.6532292366027832
$ws.Range("C49").Value = 5.2597446441650391
$ws.Range("E49").Value = 5.6242489814758301
$ws.Range("C50").Value = 5.4457979202270508
$ws.Range("E50").Value = 5.5500292778015137
$ws.Range("C51").Value = 5.4555373191833496
$ws.Range("E51").Value = 5.5456857681274414
$ws.Range("C52").Value = 5.4518184661865234
$ws.Range("D52").Value = 5.6019387245178223
$ws.Range("E52").Value = 5.4816675186157227
$ws.Range("C53").Value = 5.7958078384399414
$ws.Range("D53").Value = 5.529076099395752
$ws.Range("E53").Value = 5.3986625671386719
$ws.Range("C54").Value = 5.4417014122009277
$ws.Range("D54").Value = 5.5129237174987793
$ws.Range("E54").Value = 5.3745760917663574
$ws.Range("C55").Value = 5.569425106048584
$ws.Range("D55").Value = 5.4908051490783691
$ws.Range("E55").Value = 5.3444395065307617
$ws.Range("B56").Value = 5.3449158668518066
$ws.Range("C56").Value = 5.1834864616394043
$ws.Range("D56").Value = 5.4599113464355469
$ws.Range("E56").Value = 5.3135108947753906
$ws.Range("B57").Value = 5.1718158721923828
$ws.Range("C57").Value = 4.9846453666687012
$ws.Range("D57").Value = 5.4095525741577148
$ws.Range("E57").Value = 5.2712693214416504
$ws.Range("B58").Value = 5.1728272438049316
$ws.Range("C58").Value = 5.0429649353027344
$ws.Range("D58").Value = 5.3094196319580078
$ws.Range("E58").Value = 5.1814041137695312
$ws.Range("B59").Value = 5.2997260093688965
$ws.Range("C59").Value = 5.1745672225952148
$ws.Range("D59").Value = 5.3068537712097168
$ws.Range("E59").Value = 5.1773700714111328
$ws.Range("B60").Value = 5.3068251609802246
$ws.Range("C60").Value = 5.1771798133850098
$ws.Range("D60").Value = 5.2760229110717773
$ws.Range("E60").Value = 5.1442837715148926
$ws.Range("B61").Value = 5.0938138961791992
$ws.Range("C61").Value = 5.0716452598571777
$ws.Range("D61").Value = 5.2674112319946289
$ws.Range("E61").Value = 5.1393833160400391
$ws.Range("B62").Value = 5.0703945159912109
$ws.Range("C62").Value = 4.9870219230651855
$ws.Range("D62").Value = 5.2810678482055664
$ws.Range("E62").Value = 5.1614885330200195
$ws.Range("B63").Value = 5.5866193771362305
$ws.Range("C63").Value = 5.4053940773010254
$ws.Range("D63").Value = 5.299107551574707
$ws.Range("E63").Value = 5.1812424659729004
$ws.Range("B64").Value = 5.4372677803039551
$ws.Range("C64").Value = 5.2716474533081055
$ws.Range("D64").Value = 5.2989840507507324
$ws.Range("E64").Value = 5.1825776100158691
$ws.Range("A58").Value = 45536
$ws.Range("A59").Value = 45566
$ws.Range("A60").Value = 45597
$ws.Range("A61").Value = 45627
$ws.Range("A62").Value = 45658
$ws.Range("A63").Value = 45689
$ws.Range("A64").Value = 45717
